$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 28657.5
$ws.Range("C3").Value = 76288.5
$ws.Range("C4").Value = 1939273.94
$ws.Range("C8").Value = 42501.96
$ws.Range("C10").Value = 131397.38
$ws.Range("C11").Value = 3982.96
$ws.Range("C12").Value = 106462.84
$ws.Range("C13").Value = 3790.84
$ws.Range("C14").Value = 2178
$ws.Range("C15").Value = 7030
$ws.Range("C16").Value = 11551.69
$ws.Range("C18").Value = 19984
$ws.Range("C19").Value = 15001
$ws.Range("C20").Value = 45878.86
$ws.Range("C21").Value = 75.06999999999999
$ws.Range("C22").Value = 13518.67
$ws.Range("C23").Value = 513.62
$ws.Range("C24").Value = 3019.44
$ws.Range("C25").Value = 13.04
$ws.Range("C26").Value = 1592
$ws.Range("C27").Value = 1592
$ws.Range("C28").Value = 124387.95
$ws.Range("C30").Value = 114963.25
$ws.Range("C32").Value = 39096
$ws.Range("C33").Value = 48317
$ws.Range("C34").Value = 10514.3
$ws.Range("C35").Value = 291.07
$ws.Range("C36").Value = 1178.13
$ws.Range("C38").Value = 344446.52
$ws.Range("C42").Value = 146066.5

$ws.Range("A44").Value = "(SUMA) - E. Gómez - I. Díaz"
$ws.Range("C44").Value = 977.22

$ws.Range("C45").Value = 2533
$ws.Range("C46").Value = 2533
$ws.Range("C47").Value = 20018.16
$ws.Range("C49").Value = 102078.5
$ws.Range("C50").Value = 161181.5
$ws.Range("C51").Value = 3056156.81
$ws.Range("C52").Value = 277005.8
